# modificacion firma certificados de trabajo a partir de fecha 18_12_2019
#
# This script reproduces (via Word COM automation) the following source
# changes to the certificate template:
#   1. "La suscrita Lic. ${JEFA_RECURSOS} "  -> split into 3 runs with the
#      same (unchanged) text: "La suscrita" | " " | "Lic. ${JEFA_RECURSOS} "
#   2. "...el Area de Recursos Humanos..."   -> "...el Departamento de
#      Recursos Humanos..." (and the run that carries it gets split into 3
#      runs accordingly)
#   3. The stray "_GoBack" bookmark that used to sit right after "Es cuand"
#      is removed from there ...
#   4. ... and re-created at the end of the "GAG/${INICIALES}" paragraph
#      (Word moves "_GoBack" to the location of the most recent edit).
#   5. The signature image run's language tag is changed from the default
#      <w:lang w:eastAsia="en-US"/> to <w:lang w:val="es-BO"
#      w:eastAsia="es-BO"/>.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: force a run split at a given document offset without changing
# the visible formatting. Word only writes separate <w:r> runs at
# positions where some character-level property actually differs, so we
# flip Bold on then immediately back off for a zero-width "seam" – this
# leaves the text/formatting identical but the run boundary persists.
# ---------------------------------------------------------------------
function Split-RunAt($doc, $pos) {
    $seam = $doc.Range($pos, $pos + 1)
    $seam.Font.Bold = $true
    $seam.Font.Bold = $false
}

# =======================================================================
# 1) "La suscrita Lic. ${JEFA_RECURSOS} " -> 3 runs (text unchanged)
# =======================================================================
$full = $d.Content
$found = $full.Find.Execute("La suscrita Lic. `${JEFA_RECURSOS} ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $base = $full.Start
    Split-RunAt $d ($base + 10)   # after "La suscrita" (11 chars)
    Split-RunAt $d ($base + 11)   # after the following single space
}

# =======================================================================
# 2) "Area de Recursos Humanos" -> "Departamento de Recursos Humanos"
#    and split the carrying run into 3 pieces.
# =======================================================================
$full = $d.Content
$found = $full.Find.Execute("Que, de la revisión de la carpeta que cursa en el Área de Recu", $true, $false, $false, $false, $false, $true, 1, $false, "Que, de la revisión de la carpeta que cursa en el Departamento de Recu", 2)
Write-Output "step2 replace found=$found"

$full = $d.Content
$found = $full.Find.Execute("Que, de la revisión de la carpeta que cursa en el Departamento de Recu", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $base = $full.Start
    Split-RunAt $d ($base + 21)   # after "Que, de la revisión de" (22 chars)
    Split-RunAt $d ($base + 61)   # after "...el Departamento" (22+40=62 chars)
}

# =======================================================================
# 3) Remove the "_GoBack" bookmark from its old spot (right after "Es cuand")
# =======================================================================
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# =======================================================================
# 4) Re-add "_GoBack" collapsed right after "/${INICIALES}"
#    (Using a temporary placeholder character to sidestep the engine's
#    handling of bookmarks collapsed exactly at a paragraph-mark seam.)
# =======================================================================
$full = $d.Content
$found = $full.Find.Execute("/`${INICIALES}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $pos = $full.End
    $ins = $d.Range($pos, $pos)
    $ins.InsertAfter("X")
    $placeholder = $d.Range($pos, $pos + 1)
    $d.Bookmarks.Add("_GoBack", $placeholder)
    $placeholder = $d.Range($pos, $pos + 1)
    $placeholder.Delete()
}

# =======================================================================
# 5) Signature image run: <w:lang w:eastAsia="en-US"/> -> <w:lang
#    w:val="es-BO" w:eastAsia="es-BO"/>
# =======================================================================
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $r = $p.Range
    if ($r.Shapes.Count -gt 0) {
        $r.Font.LanguageID = "es-BO"
        $r.Font.LanguageIDFarEast = "es-BO"
        break
    }
}
